# KIBON-250: InstitutionStammdaten: Felder oeffnungstage und oeffnungsstunden entfernen
#
# The "Kanton.xlsx" report template had a column showing the institution's
# "Öffnungstage" (opening days) with its "{oeffnungstage}" merge placeholder.
# That field is being removed from the model, so the column that renders it
# is deleted from the report sheet as well.
#
# In the original workbook that column is "S" (header "Öffnungstage" in S7,
# merged S7:S8, and the placeholder "{oeffnungstage}" in S9). The very last
# column, "T", only held the "{repeatKantonRow}" marker used by the report
# engine to know where to repeat rows; once S is removed, T's content simply
# slides left into the new column S - exactly what Excel's own
# "delete entire column" does to every cell/format/merge/column-width to the
# right of the deleted column.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete the whole "Öffnungstage" column (S). This removes the header cell,
# the S7:S8 merge, the "{oeffnungstage}" placeholder, and shifts every
# column to its right (T, which held "{repeatKantonRow}") one slot to the
# left.
$ws.Range("S1").EntireColumn.Delete() | Out-Null

# Put the active selection back onto the (now shifted) last used row, as in
# the saved workbook.
$ws.Range("S17").Select() | Out-Null
